# Update column F ("dSF") values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$updates = @{
    2  = -2
    3  = -6
    4  = -2
    5  = 2
    6  = -4
    7  = 4
    8  = -6
    10 = 7
    11 = 3
    12 = 4
    13 = -4
    14 = 3
    15 = 2
    16 = 3
    18 = -5
    19 = -4
    20 = -2
    21 = -2
    22 = -3
    23 = 1
    24 = -1
    25 = -1
    26 = 1
    27 = 3
    28 = -5
    31 = 1
    32 = -2
    33 = -2
    34 = -2
    35 = -2
    36 = -1
    37 = 3
    38 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
